$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the underlying expiration label text (was "OEW4 Feb 24, 2017")
$ws.Range("B2").Value = "OEW Feb 28, 2017"

# Update the expiration date value (was 2/24/2017 -> now 2/28/2017)
$ws.Range("C2").Value = "2/28/2017"

# Update the current selection to match the new editing position
$ws.Range("F14").Select()
